# Update the player/position/team table on Sheet1 to reflect the reshuffled
# roster data from the upload. Row 1 (headers) and formatting are unchanged;
# only the B (Pozisyon) and C (Takım) values for rows 2-19 change, and the
# row data itself is reassigned to different players/teams as per the new
# upload.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @(2,  "Jose Alvarado",            "PG",         "New Orleans Pelicans"),
    @(3,  "RJ Barrett",                "SG,SF,PF",   "Toronto Raptors"),
    @(4,  "CJ McCollum",               "PG,SG",      "New Orleans Pelicans"),
    @(5,  "Quentin Grimes",            "SG,SF",      "Dallas Mavericks"),
    @(6,  "Zach LaVine",               "SG,SF",      "Chicago Bulls"),
    @(7,  "Nikola Jovic",              "PF,C",       "Miami Heat"),
    @(8,  "Guerschon Yabusele",        "PF,C",       "Philadelphia 76ers"),
    @(9,  "Lauri Markkanen",           "SF,PF",      "Utah Jazz"),
    @(10, "Shai Gilgeous-Alexander",   "PG,SG",      "Oklahoma City Thunder"),
    @(11, "Toumani Camara",            "SF,PF",      "Portland Trail Blazers"),
    @(12, "Andrew Nembhard",           "PG,SG",      "Indiana Pacers"),
    @(13, "Jalen Williams",            "SG,SF,PF,C", "Oklahoma City Thunder"),
    @(14, "Tobias Harris",             "SF,PF",      "Detroit Pistons"),
    @(15, "Jordan Poole",              "PG,SG",      "Washington Wizards"),
    @(16, "Kyrie Irving",              "PG,SG",      "Dallas Mavericks"),
    @(17, "John Collins",              "PF,C",       "Utah Jazz"),
    @(18, "Joel Embiid",               "C",          "Philadelphia 76ers"),
    @(19, "Jimmy Butler",              "SF,PF",      "Miami Heat")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}
